$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 500657
$ws.Range("J3").Value = 500657
$ws.Range("L3").Value = 500657
$ws.Range("N3").Value = -500885

# Row 11
$ws.Range("H11").Value = 11.75
$ws.Range("I11").Value = 11.75
$ws.Range("K11").Value = 11.75
$ws.Range("M11").Value = 128.25

# Row 17
$ws.Range("H17").Value = 1744.9231
$ws.Range("J17").Value = 1744.9231
$ws.Range("L17").Value = 5234.7693
$ws.Range("N17").Value = -5570.7693

# Row 31
$ws.Range("H31").Value = 111111450
$ws.Range("I31").Value = 111111450
$ws.Range("K31").Value = 333334350
$ws.Range("M31").Value = -333334120

# Row 32
$ws.Range("H32").Value = 2325.9375
$ws.Range("I32").Value = 1749.6
$ws.Range("J32").Value = 2587.9092
$ws.Range("K32").Value = 1749.6
$ws.Range("L32").Value = 2587.9092
$ws.Range("M32").Value = -1423.6
$ws.Range("N32").Value = -3239.9092

# Row 38
$ws.Range("H38").Value = 365.1111
$ws.Range("I38").Value = 133.25
$ws.Range("K38").Value = 399.75
$ws.Range("M38").Value = -27.75

# Row 39
$ws.Range("H39").Value = 631.9167
$ws.Range("I39").Value = 680.2727
$ws.Range("K39").Value = 2040.8181
$ws.Range("M39").Value = -1744.8181

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 74
$ws.Range("H74").Value = 107343.266
$ws.Range("I74").Value = 107343.266
$ws.Range("K74").Value = 107343.266
$ws.Range("M74").Value = -106407.266

# Row 77
$ws.Range("H77").Value = 107343.266
$ws.Range("I77").Value = 107343.266
$ws.Range("K77").Value = 536716.3300000001
$ws.Range("M77").Value = -532036.3300000001

# Row 102
$ws.Range("H102").Value = 500657
$ws.Range("J102").Value = 500657
$ws.Range("L102").Value = 500657
$ws.Range("N102").Value = -507147

# Row 107
$ws.Range("H107").Value = 545.5
$ws.Range("I107").Value = 571.5
$ws.Range("J107").Value = 519.5
$ws.Range("K107").Value = 571.5
$ws.Range("L107").Value = 519.5
$ws.Range("M107").Value = 1348.5
$ws.Range("N107").Value = -4359.5

# Row 138
$ws.Range("H138").Value = 4215.8823
$ws.Range("I138").Value = 2096.1333
$ws.Range("J138").Value = 5889.3687
$ws.Range("K138").Value = 6288.3999
$ws.Range("L138").Value = 17668.1061
$ws.Range("M138").Value = -1148.3999
$ws.Range("N138").Value = -27948.1061

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2609.0715
$ws.Range("I32").Value = 2690.1428
$ws.Range("J32").Value = 2203.7144
$ws.Range("K32").Value = 2690.1428
$ws.Range("L32").Value = 2203.7144
$ws.Range("M32").Value = -2403.1428
$ws.Range("N32").Value = -2777.7144

# Row 61
$ws.Range("H61").Value = 3231.875
$ws.Range("I61").Value = 3231.875
$ws.Range("K61").Value = 3231.875
$ws.Range("M61").Value = -3019.875

# Row 136
$ws.Range("H136").Value = 3231.875
$ws.Range("I136").Value = 3231.875
$ws.Range("K136").Value = 9695.625
$ws.Range("M136").Value = -7145.625

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3058.5
$ws.Range("I86").Value = 3212.125
$ws.Range("J86").Value = 2444
$ws.Range("K86").Value = 3212.125
$ws.Range("L86").Value = 2444
$ws.Range("M86").Value = -2089.125
$ws.Range("N86").Value = -4690

# Row 89
$ws.Range("H89").Value = 3058.5
$ws.Range("I89").Value = 3212.125
$ws.Range("J89").Value = 2444
$ws.Range("K89").Value = 16060.625
$ws.Range("L89").Value = 12220
$ws.Range("M89").Value = -10444.625
$ws.Range("N89").Value = -23452

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1135.8
$ws.Range("I107").Value = 1054.5714
$ws.Range("J107").Value = 1325.3334
$ws.Range("K107").Value = 1054.5714
$ws.Range("L107").Value = 1325.3334
$ws.Range("M107").Value = 865.4286
$ws.Range("N107").Value = -5165.3334

# Row 141
$ws.Range("H141").Value = 41994.25
$ws.Range("J141").Value = 40996
$ws.Range("L141").Value = 40996
$ws.Range("N141").Value = -51356

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 931.7857
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 965
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 2895
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = -3119

# Row 113
$ws.Range("H113").Value = 1946.5834
$ws.Range("J113").Value = 2066
$ws.Range("L113").Value = 6198
$ws.Range("N113").Value = -10538

# Row 128
$ws.Range("H128").Value = 624073.5
$ws.Range("I128").Value = 624073.5
$ws.Range("K128").Value = 1872220.5
$ws.Range("M128").Value = -1867240.5

# Row 132
$ws.Range("H132").Value = 3575
$ws.Range("I132").Value = 704
$ws.Range("J132").Value = 4532
$ws.Range("K132").Value = 6336
$ws.Range("L132").Value = 40788
$ws.Range("M132").Value = -3806
$ws.Range("N132").Value = -45848

# Row 135
$ws.Range("H135").Value = 931.7857
$ws.Range("I135").Value = 500
$ws.Range("J135").Value = 965
$ws.Range("K135").Value = 4500
$ws.Range("L135").Value = 8685
$ws.Range("M135").Value = -1965
$ws.Range("N135").Value = -13755

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8605.799999999999
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 8605.799999999999
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 8605.799999999999
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -8829.799999999999

# Row 126
$ws.Range("H126").Value = 8605.799999999999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8605.799999999999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 25817.4
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -30757.4

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 133.28572
$ws.Range("I18").Value = 133.28572
$ws.Range("K18").Value = 133.28572
$ws.Range("M18").Value = 39.71428

# Row 20
$ws.Range("H20").Value = 3505.5
$ws.Range("J20").Value = 3505.5
$ws.Range("L20").Value = 3505.5
$ws.Range("N20").Value = -3985.5

# Row 122
$ws.Range("H122").Value = 2471.3845
$ws.Range("I122").Value = 2011.7273
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 6035.1819
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -3585.1819
$ws.Range("N122").Value = -19898.5
